$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-13 (columns A-T), reflecting the updated NATMI LR-pair stats
# following additional replicate (Dr Hou advice).
$data = @(
    @("ECs", "Il16", "Kcnj15", "ECs", 3.0, 1.0, 2.528209, 7.584626999999999, 0.0779745999616348, 0.0779745999616348, 3.0, 1.0, 1.976664666666667, 5.929994, 0.8737289477891174, 0.8737289477891172, 4.997421400248666, 44.97679260223799, 0.06812866517875653, 0.06812866517875651),
    @("ECs", "Il16", "Kcnj15", "FAPs", 3.0, 1.0, 2.528209, 7.584626999999999, 0.0779745999616348, 0.0779745999616348, 1.0, 0.3333333333333333, 0.285667, 0.857001, 0.1262710522108828, 0.1262710522108827, 0.722225880403, 6.500032923627, 0.009845934782878285, 0.009845934782878282),
    @("FAPs", "Il16", "Kcnj15", "ECs", 3.0, 1.0, 1.344088666666667, 4.032266, 0.04145415829794943, 0.04145415829794944, 3.0, 1.0, 1.976664666666667, 5.929994, 0.8737289477891174, 0.8737289477891172, 2.656812576267111, 23.911313186404, 0.03621969811115087, 0.03621969811115087),
    @("FAPs", "Il16", "Kcnj15", "FAPs", 3.0, 1.0, 1.344088666666667, 4.032266, 0.04145415829794943, 0.04145415829794944, 1.0, 0.3333333333333333, 0.285667, 0.857001, 0.1262710522108828, 0.1262710522108827, 0.3839617771406667, 3.455655994266, 0.005234460186798571, 0.005234460186798571),
    @("M1", "Il16", "Kcnj15", "ECs", 3.0, 1.0, 5.501111666666668, 16.503335, 0.1696643677609784, 0.1696643677609785, 3.0, 1.0, 1.976664666666667, 5.929994, 0.8737289477891174, 0.8737289477891172, 10.87385305888778, 97.86467752999002, 0.1482406695211055, 0.1482406695211055),
    @("M1", "Il16", "Kcnj15", "FAPs", 3.0, 1.0, 5.501111666666668, 16.503335, 0.1696643677609784, 0.1696643677609785, 1.0, 0.3333333333333333, 0.285667, 0.857001, 0.1262710522108828, 0.1262710522108827, 1.571486066481667, 14.143374598335, 0.02142369823987292, 0.02142369823987292),
    @("M2", "Il16", "Kcnj15", "ECs", 3.0, 1.0, 9.507065666666666, 28.521197, 0.293215332342906, 0.2932153323429061, 3.0, 1.0, 1.976664666666667, 5.929994, 0.8737289477891174, 0.8737289477891172, 18.79228078697977, 169.130527082818, 0.2561907238036036, 0.2561907238036036),
    @("M2", "Il16", "Kcnj15", "FAPs", 3.0, 1.0, 9.507065666666666, 28.521197, 0.293215332342906, 0.2932153323429061, 1.0, 0.3333333333333333, 0.285667, 0.857001, 0.1262710522108828, 0.1262710522108827, 2.715854927799667, 24.442694350197, 0.03702460853930242, 0.03702460853930242),
    @("Neutro", "Il16", "Kcnj15", "ECs", 3.0, 1.0, 13.22109733333333, 39.663292, 0.4077628770487342, 0.4077628770487343, 3.0, 1.0, 1.976664666666667, 5.929994, 0.8737289477891174, 0.8737289477891172, 26.13367595336089, 235.203083580248, 0.3562742295112538, 0.3562742295112538),
    @("Neutro", "Il16", "Kcnj15", "FAPs", 3.0, 1.0, 13.22109733333333, 39.663292, 0.4077628770487342, 0.4077628770487343, 1.0, 0.3333333333333333, 0.285667, 0.857001, 0.1262710522108828, 0.1262710522108827, 3.776831211921333, 33.991480907292, 0.05148864753748048, 0.05148864753748048),
    @("sCs", "Il16", "Kcnj15", "ECs", 3.0, 1.0, 0.321922, 0.9657660000000001, 0.009928664587796896, 0.009928664587796898, 3.0, 1.0, 1.976664666666667, 5.929994, 0.8737289477891174, 0.8737289477891172, 0.6363318428226667, 5.726986585404, 0.008674961663246853, 0.008674961663246853),
    @("sCs", "Il16", "Kcnj15", "FAPs", 3.0, 1.0, 0.321922, 0.9657660000000001, 0.009928664587796896, 0.009928664587796898, 1.0, 0.3333333333333333, 0.285667, 0.857001, 0.1262710522108828, 0.1262710522108827, 0.09196249197400001, 0.8276624277660001, 0.001253702924550045, 0.001253702924550045),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $colNum = $j + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}
